$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 1: Collapse the three detailed CORE COMPETENCIES paragraphs
# into a single summary paragraph.
# ---------------------------------------------------------------------
$bullet = [char]0x2022
$p6 = $d.Paragraphs.Item(6)
$p6.Range.Text = "Data Visualization & Design " + $bullet + " Geospatial Analysis & Mapping " + $bullet + " Technical Visualization"

# The two now-redundant detail paragraphs (still at indices 7 and 8,
# since only text inside paragraph 6 changed) are removed entirely.
$p7 = $d.Paragraphs.Item(7)
$p8 = $d.Paragraphs.Item(8)
$rangeToDelete = $d.Range($p7.Range.Start, $p8.Range.End)
$rangeToDelete.Delete()

# ---------------------------------------------------------------------
# Change 2: Add a new "TECHNICAL SKILLS" section at the end of the
# document, right before the closing "For a more detailed..." line.
# ---------------------------------------------------------------------
$total = $d.Paragraphs.Count
$pAnchor = $d.Paragraphs.Item($total - 1)

$pAnchor.Range.InsertParagraphAfter()
$pHeading = $d.Paragraphs.Item($total)
$pHeading.Range.Text = "TECHNICAL SKILLS"

$pHeading.Range.InsertParagraphAfter()
$pSkill1 = $d.Paragraphs.Item($total + 1)
$pSkill1.Range.Text = "DATA VISUALIZATION & DESIGN Interactive Dashboards; Statistical Visualization; Geospatial Mapping; Choropleth Design"

$pSkill1.Range.InsertParagraphAfter()
$pSkill2 = $d.Paragraphs.Item($total + 2)
$pSkill2.Range.Text = "GEOSPATIAL ANALYSIS & MAPPING Spatial Analysis; Mapping Technologies; Web Mapping; Spatial Data Processing"

$pSkill2.Range.InsertParagraphAfter()
$pSkill3 = $d.Paragraphs.Item($total + 3)
$pSkill3.Range.Text = "TECHNICAL VISUALIZATION Programming; Database Integration; Web Technologies; Statistical Computing"

# Apply the Heading 2 style only to the section title, after all the
# sibling paragraphs have been created (so they don't inherit it).
$pHeading.Style = $d.Styles.Item("Heading 2")
